# Update power-flow results (pl_mw) for the "380 kV" case: rows 2-25 (A=0..23),
# columns B, C, E, F, G, H, J, K get new simulation values. Columns A, D, I, L, M, N, O
# are left untouched (they stay 0 / unchanged index values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update range B2:C25
$data_B2 = New-Object 'object[,]' 24,2
$data_B2[0,0] = 0.435511965976275
$data_B2[0,1] = 0.04311331021581566
$data_B2[1,0] = 0.3968284938217153
$data_B2[1,1] = 0.03769916303504317
$data_B2[2,0] = 0.3732193453386117
$data_B2[2,1] = 0.03436194866455367
$data_B2[3,0] = 0.363634512743829
$data_B2[3,1] = 0.03299877983344857
$data_B2[4,0] = 0.3620451446306561
$data_B2[4,1] = 0.03277223236396765
$data_B2[5,0] = 0.3730899345305829
$data_B2[5,1] = 0.03434357753312156
$data_B2[6,0] = 0.4221444098546954
$data_B2[6,1] = 0.0412492000925937
$data_B2[7,0] = 0.5194680708577266
$data_B2[7,1] = 0.05468858933524245
$data_B2[8,0] = 0.5916621584337349
$data_B2[8,1] = 0.06450089864155473
$data_B2[9,0] = 0.6246562168807372
$data_B2[9,1] = 0.06895166521654517
$data_B2[10,0] = 0.6371720837881014
$data_B2[10,1] = 0.07063520101330312
$data_B2[11,0] = 0.6344756046302109
$data_B2[11,1] = 0.0702727050235552
$data_B2[12,0] = 0.6256854704102466
$data_B2[12,1] = 0.06909020835232127
$data_B2[13,0] = 0.6203040835028162
$data_B2[13,1] = 0.06836565058769395
$data_B2[14,0] = 0.5895089717186579
$data_B2[14,1] = 0.06420976942382595
$data_B2[15,0] = 0.5706561016710623
$data_B2[15,1] = 0.06165695349099565
$data_B2[16,0] = 0.5598268057022437
$data_B2[16,1] = 0.06018742785522591
$data_B2[17,0] = 0.5561626678403115
$data_B2[17,1] = 0.05968966387159469
$data_B2[18,0] = 0.5726615355369233
$data_B2[18,1] = 0.06192883082238154
$data_B2[19,0] = 0.6282667577331154
$data_B2[19,1] = 0.06943758725326177
$data_B2[20,0] = 0.6647345679048726
$data_B2[20,1] = 0.0743340695640029
$data_B2[21,0] = 0.6452594972064105
$data_B2[21,1] = 0.07172172986224723
$data_B2[22,0] = 0.5717548493734057
$data_B2[22,1] = 0.0618059209307944
$data_B2[23,0] = 0.4930184234258377
$data_B2[23,1] = 0.05106376627652764
$ws.Range("B2:C25").Value = $data_B2

# Update range E2:H25
$data_E2 = New-Object 'object[,]' 24,4
$data_E2[0,0] = 0.4987668589820373
$data_E2[0,1] = 2.660169050887376
$data_E2[0,2] = 0.9796289912398208
$data_E2[0,3] = 1.004526499947474
$data_E2[1,0] = 0.4765312875846206
$data_E2[1,1] = 2.625888965376717
$data_E2[1,2] = 0.9817520460068039
$data_E2[1,3] = 1.010959240581656
$data_E2[2,0] = 0.4631162117635199
$data_E2[2,1] = 2.606455037019487
$data_E2[2,2] = 0.9837600347720752
$data_E2[2,3] = 1.015417111240936
$data_E2[3,0] = 0.4577091545523402
$data_E2[3,1] = 2.598940462489054
$data_E2[3,2] = 0.9847548364299996
$data_E2[3,3] = 1.017361337397233
$data_E2[4,0] = 0.4568149204895633
$data_E2[4,1] = 2.597717099873975
$data_E2[4,2] = 0.9849306668009064
$data_E2[4,3] = 1.017691876932744
$data_E2[5,0] = 0.4630430486775197
$data_E2[5,1] = 2.606352054697624
$data_E2[5,2] = 0.9837727370580041
$data_E2[5,3] = 1.015442815305605
$data_E2[6,0] = 0.4910506601162723
$data_E2[6,1] = 2.648013638095847
$data_E2[6,2] = 0.9802144520858747
$data_E2[6,3] = 1.006638941856636
$data_E2[7,0] = 0.5478657138679495
$data_E2[7,1] = 2.742572348793061
$data_E2[7,2] = 0.9788541898349479
$data_E2[7,3] = 0.9934144274795358
$data_E2[8,0] = 0.5907755003950967
$data_E2[8,1] = 2.81996916408707
$data_E2[8,2] = 0.9813199474278349
$data_E2[8,3] = 0.9861727572085357
$data_E2[9,0] = 0.6105533046137452
$data_E2[9,1] = 2.856918768715303
$data_E2[9,2] = 0.9832029390410355
$data_E2[9,3] = 0.9834181769337533
$data_E2[10,0] = 0.6180799160185586
$data_E2[10,1] = 2.871162365618801
$data_E2[10,2] = 0.9840261451501675
$data_E2[10,3] = 0.982452901869209
$data_E2[11,0] = 0.6164572701184738
$data_E2[11,1] = 2.868083544143246
$data_E2[11,2] = 0.9838439427852705
$data_E2[11,3] = 0.9826573263698748
$data_E2[12,0] = 0.6111717776230705
$data_E2[12,1] = 2.858085549143823
$data_E2[12,2] = 0.9832684530549471
$data_E2[12,3] = 0.9833372021997491
$data_E2[13,0] = 0.6079391065596269
$data_E2[13,1] = 2.85199428709592
$data_E2[13,2] = 0.9829303151011715
$data_E2[13,3] = 0.9837637873306448
$data_E2[14,0] = 0.5894881726519827
$data_E2[14,1] = 2.817589561934653
$data_E2[14,2] = 0.98121226729738
$data_E2[14,3] = 0.9863636540664515
$data_E2[15,0] = 0.5782352505167268
$data_E2[15,1] = 2.796930223002278
$data_E2[15,2] = 0.98035377515059
$data_E2[15,3] = 0.988096970631986
$data_E2[16,0] = 0.5717871385246127
$data_E2[16,1] = 2.785211356165462
$data_E2[16,2] = 0.979931609990885
$data_E2[16,3] = 0.9891447137187299
$data_E2[17,0] = 0.5696080843491984
$data_E2[17,1] = 2.781271653140237
$data_E2[17,2] = 0.9798009515873929
$data_E2[17,3] = 0.9895081773348267
$data_E2[18,0] = 0.5794306306323307
$data_E2[18,1] = 2.799112479769491
$data_E2[18,2] = 0.9804377460820746
$data_E2[18,3] = 0.9879071984634606
$data_E2[19,0] = 0.6127232448702671
$data_E2[19,1] = 2.861015365328342
$data_E2[19,2] = 0.9834344932307602
$data_E2[19,3] = 0.9831353923066501
$data_E2[20,0] = 0.6346986914017805
$data_E2[20,1] = 2.90293945319317
$data_E2[20,2] = 0.9860354439892092
$data_E2[20,3] = 0.9804704649793621
$data_E2[21,0] = 0.6229501114832914
$data_E2[21,1] = 2.880429160105479
$data_E2[21,2] = 0.9845882616789368
$data_E2[21,3] = 0.9818511988484033
$data_E2[22,0] = 0.5788901328129725
$data_E2[22,1] = 2.79812538795602
$data_E2[22,2] = 0.9803995604800235
$data_E2[22,3] = 0.9879928348682796
$data_E2[23,0] = 0.5322917124527038
$data_E2[23,1] = 2.715606221676808
$data_E2[23,2] = 0.9786165938650697
$data_E2[23,3] = 0.9965583223636827
$ws.Range("E2:H25").Value = $data_E2

# Update range J2:K25
$data_J2 = New-Object 'object[,]' 24,2
$data_J2[0,0] = 0.07823170415101899
$data_J2[0,1] = 0.4016520619536834
$data_J2[1,0] = 0.07861575489026862
$data_J2[1,1] = 0.3601842550331469
$data_J2[2,0] = 0.07890149843764505
$data_J2[2,1] = 0.3348104708113624
$data_J2[3,0] = 0.07903046092994259
$data_J2[3,1] = 0.3244926836564161
$data_J2[4,0] = 0.07905263013456576
$data_J2[4,1] = 0.3227807721396516
$data_J2[5,0] = 0.07890318702793309
$data_J2[5,1] = 0.3346712312169018
$data_J2[6,0] = 0.0783537385834947
$data_J2[6,1] = 0.3873358844625727
$data_J2[7,0] = 0.07767419697007227
$data_J2[7,1] = 0.4913030581764133
$data_J2[8,0] = 0.07741997346184704
$data_J2[8,1] = 0.5681136478917779
$data_J2[9,0] = 0.07735806280309987
$data_J2[9,1] = 0.6031507609774565
$data_J2[10,0] = 0.07734238878969535
$data_J2[10,1] = 0.616432100014606
$data_J2[11,0] = 0.07734541820147811
$data_J2[11,1] = 0.6135711267405384
$data_J2[12,0] = 0.07735661731934584
$data_J2[12,1] = 0.6042431544039175
$data_J2[13,0] = 0.07736449028954695
$data_J2[13,1] = 0.5985312600826092
$data_J2[14,0] = 0.07742510445647355
$data_J2[14,1] = 0.5658257963272604
$data_J2[15,0] = 0.07747608306609521
$data_J2[15,1] = 0.5457863913582912
$data_J2[16,0] = 0.07751046000928952
$data_J2[16,1] = 0.5342692646569844
$data_J2[17,0] = 0.07752296641522349
$data_J2[17,1] = 0.5303713169967068
$data_J2[18,0] = 0.07747013278123305
$data_J2[18,1] = 0.5479186893086307
$data_J2[19,0] = 0.07735311664286115
$data_J2[19,1] = 0.6069826399601652
$data_J2[20,0] = 0.07732195060952662
$data_J2[20,1] = 0.6456632067618955
$data_J2[21,0] = 0.07733442434191673
$data_J2[21,1] = 0.6250115091367547
$data_J2[22,0] = 0.07747280712029081
$data_J2[22,1] = 0.5469546655748161
$data_J2[23,0] = 0.07781517772532354
$data_J2[23,1] = 0.4631027105915848
$ws.Range("J2:K25").Value = $data_J2
